$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: shift headers one column left, dropping the old "Name" column ---
$ws.Range("B1").Value = "Potential [V]"
$ws.Range("C1").Value = "CDL [F]"
$ws.Range("D1").Value = "b [F/mV/s]"
$ws.Range("E1").Clear()

# --- New data row 2 ---
# A2: numeric 0, carrying the same (bold/bordered/centered) header style as row 1
$ws.Range("A2").Value = 0
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# B2:D2: values that look numeric but must be stored as literal text, with no
# special number-format/style attached (match via formula + paste-as-values so
# Excel doesn't auto-apply a "text" number format/style to the cell).
$ws.Range("B2").Formula = "=""0.25356765"""
$ws.Range("C2").Formula = "=""4.06e-08"""
$ws.Range("D2").Formula = "=""1.08e-06"""
$ws.Range("B2:D2").Copy()
$ws.Range("B2:D2").PasteSpecial(-4163)

$excel.CutCopyMode = 0
